# Austin Wordle Scores - add a new day's data, drop the oldest day, and
# correct a few recent scores (per commit "Adding markdown and cleaning data").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks one row per Wordle day (col A = day number, col B = score).
# The oldest tracked day (318) is being retired: delete its row so every
# later day shifts up by one.
$ws.Rows(2).Delete()

# A handful of recent scores (days 342-345, now sitting at rows 25-28) had
# the wrong value recorded - correct them.
$ws.Range("B25").Value = 2
$ws.Range("B26").Value = 4
$ws.Range("B27").Value = 3
$ws.Range("B28").Value = 4

# Append the newly played days (346-359) at the bottom of the table.
$newDays = @(
    @(346, 3),
    @(347, 5),
    @(348, 5),
    @(349, 4),
    @(350, 6),
    @(351, 4),
    @(352, 4),
    @(353, 4),
    @(354, 3),
    @(355, 4),
    @(356, 4),
    @(357, 4),
    @(358, 3),
    @(359, 4)
)

$row = 29
foreach ($entry in $newDays) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

# Update the view so it is scrolled down to the new bottom-of-data rows and
# the last-edited cell (B28) is the active selection.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B28").Select()
